$d = $word.ActiveDocument

# --- Fix bookmark id for 'introducción-histórica' is handled automatically by the engine's
# bookmark renumbering on save; we just need the name+nesting to stay correct. ---

# 1) Insert a placeholder paragraph right after the 'Las reglas prescriptivas vs. ...' heading
#    (paragraph 6) and seed a fresh decimal numbering definition (numId 1003) on it so that
#    word/numbering.xml gains a matching abstractNum/num pair before we overwrite its content.
$lastHeadingPara = $d.Paragraphs(6).Range
$lastHeadingPara.InsertParagraphAfter()
$seed = $d.Paragraphs(7).Range
$seed.Text = "seed"
$seed.ListFormat.ApplyNumberDefault()

# 2) Replace that seeded paragraph with the full block of new content (8 paragraphs).
$block1Target = $d.Paragraphs(7).Range
$block1Xml = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:rPr><w:bCs/><w:b/></w:rPr><w:t xml:space="preserve">Las reglas</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/><w:bCs/><w:b/></w:rPr><w:t xml:space="preserve">prescriptivas:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Indican cómo</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space="preserve">debería</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">ser el lenguaje, cómo</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space="preserve">deberían</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">usarla los hablantes y qué funciones y usos</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space="preserve">deberían</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">tener los elementos que lo componen.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Por ejemplo, en la oración en (1), la lingüística prescriptiva diría que se debería usar</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space="preserve">lo</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">en vez de</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space="preserve">le</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">y</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space="preserve">había</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">en (2).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1003"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">María ha visto a Juan. Yo también</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space="preserve">le</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">he visto.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1003"/></w:numPr></w:pPr><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space="preserve">Habían</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">tres coches aparcados.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:rPr><w:bCs/><w:b/></w:rPr><w:t xml:space="preserve">Las reglas</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/><w:bCs/><w:b/></w:rPr><w:t xml:space="preserve">descriptivas:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Indican qué construcciones se usan en realidad y las circunstancias en qué se usan.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Por ejemplo, la lingüística descriptiva diría que se usa</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space="preserve">le</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">en la zona de norte de España (y en otras partes de España), cuando el referente es masculino y animado. Y que se usa</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">habían</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">para reflejar que el sujeto de la oración es plural, como con</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">Son las tres/Es la una.</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:rPr><w:bCs/><w:b/></w:rPr><w:t xml:space="preserve">Una pregunta:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">¿Es la lingüística tradicional prescriptiva o descriptiva?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:rPr><w:bCs/><w:b/></w:rPr><w:t xml:space="preserve">Importante:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">La lingüística moderna es una lingüística descriptiva y es la lingüística que vamos a hacer nosotros</w:t></w:r></w:p></pkg:xmlData>'
$block1Target.InsertXML($block1Xml)

# 3) Append the new 'Hacía la lingüística moderna' Heading2 paragraph after the block.
$lastNewPara = $d.Paragraphs(14).Range
$lastNewPara.InsertParagraphAfter()
$headingTarget = $d.Paragraphs(15).Range
$headingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t xml:space="preserve">Hacía la lingüística moderna</w:t></w:r></w:p>'
$headingTarget.InsertXML($headingXml)

# 4) Wrap that heading paragraph with a new bookmark named 'hacía-la-lingüística-moderna'.
$headingPara = $d.Paragraphs(15)
$bmRange = $d.Range($headingPara.Range.Start, $headingPara.Range.End)
$d.Bookmarks.Add('hacía-la-lingüística-moderna', $bmRange)

Write-Host "Final paragraph count:" $d.Paragraphs.Count
